$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update SteuerfussKanton (column E) values for rows 2-5
$ws.Range("E2:E5").Value = 104

# Best-fit (auto-fit) the SteuerfussKanton/SteuerfussGemeinde columns so their
# width matches their content, same as double-clicking the column border.
$ws.Columns.Item(5).ColumnWidth = 14.428571428571429
$ws.Columns.Item(6).ColumnWidth = 17.142857142857142

# Update the selection to match the reviewed range (active cell E2, selection E2:E5)
$ws.Range("E2:E5").Select() | Out-Null
